$d = $word.ActiveDocument
$W_NS = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Locate the "person?.name = ..." paragraph and strip its stray
#     paragraph-/run-level formatting (pPr + rPr lang="en-US"), keeping the text as-is.
$rngNameQ = $d.Content
$null = $rngNameQ.Find.Execute("person?.name = nullable operator auf name, existiert name?", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pNameQ = $rngNameQ.Paragraphs(1)
$xmlNameQ = "<w:p $W_NS><w:r><w:t>person?.name = nullable operator auf name, existiert name?</w:t></w:r></w:p>"
$null = $pNameQ.Range.InsertXML($xmlNameQ)

# --- Locate the "person | json = person als jso" + "n ausgeben" paragraph (which also
#     carries the _GoBack bookmark) and replace it with the corrected single-run paragraph,
#     then append the new notes that were added: an empty separator paragraph, the
#     "#erstellt in subfolder" paragraph (now carrying the relocated bookmark), and the
#     "ng g c accounts/acountlist" paragraph.
$rngJson = $d.Content
$null = $rngJson.Find.Execute("person | json = person als jso", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pJson = $rngJson.Paragraphs(1)
$xmlJson = "<w:p $W_NS><w:r><w:t>person | json = person als json ausgeben</w:t></w:r></w:p>"
$xmlJson += "<w:p $W_NS/>"
$xmlJson += "<w:p $W_NS><w:r><w:t>#erstellt in subfolder</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"
$xmlJson += "<w:p $W_NS><w:r><w:t>ng g c accounts/acountlist</w:t></w:r></w:p>"
$null = $pJson.Range.InsertXML($xmlJson)

# --- The final (now trailing) empty paragraph also carried stray "en-US" paragraph
#     formatting; normalise it to a bare empty paragraph.
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$xmlLast = "<w:p $W_NS/>"
$null = $pLast.Range.InsertXML($xmlLast)
